# TC_ACC_11/Default.xlsx: sample test script cleanup.
# The "Global" sheet used to hold a sample login row (USERNAME/PASSWORD
# header with a numeric id + "Qatar@2021" password underneath). The new
# scripts only need a tiny "A"/"B" header row, so the old credential row
# is removed and the headers are renamed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Global")

# Rename the header cells: B first so the shared-string table order comes
# out as [0]="B", [1]="A" (matches the authored workbook).
$ws.Range("B1").Value = "B"
$ws.Range("A1").Value = "A"

# Drop the old sample credential row (id 88996 / "Qatar@2021") completely.
$ws.Rows.Item(2).Delete()

# The remaining header cells no longer need their old best-fit widths.
$ws.Range("A1:B1").ColumnWidth = 8.6

# Move the sheet's selection to C3 ...
$ws.Range("C3").Select()

# ... then restore "Action1" as the active/selected sheet/tab, since that
# was the active sheet before this edit and should stay that way.
$wb.Worksheets.Item("Action1").Select()
